$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove stale contents from the old layout (A1:E1) before laying out the new table
$ws.Range("A1:E1").ClearContents()

# Header row
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "pierdoli"
$ws.Range("D1").Value = "smiedzi"

# Row 2
$ws.Range("B2").Value = "debil@gmail.com"
$ws.Range("C2").Value = "no debil no"
$ws.Range("D2").Value = "N/A"

# Row 3
$ws.Range("B3").Value = "idiota@gmail.com"
$ws.Range("C3").Value = "N/A"
$ws.Range("D3").Value = "no idiota no"

# A2/A3 hold "1"/"2" but as *text* (matches the shared-string-backed source
# workbook) rather than numbers. A plain .Value assignment gets re-parsed as
# a number, so round-trip them through a text formula + paste-values instead
# of flipping NumberFormat to "@" (which would permanently register a new
# cell style even after the helper cell is cleared).
$helper = $ws.Range("Z100")

$helper.Formula = "=""1"""
$helper.Copy()
$ws.Range("A2").PasteSpecial(-4163)

$helper.Formula = "=""2"""
$helper.Copy()
$ws.Range("A3").PasteSpecial(-4163)

$helper.ClearContents()
